$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("F2").Value = 6662.05
$wsSummary.Range("A3").Value = 31.3
$wsSummary.Range("E3").Value = 31.3
$wsSummary.Range("F3").Value = 18.91

# ---------------------------------------------------------------------------
# Sheet: Repayment schedule
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Row 9
$wsSchedule.Range("B9").Value = 0
$wsSchedule.Range("C9").Value = 42093
$wsSchedule.Range("F9").Value = 835.12
$wsSchedule.Range("G9").Value = 4177.2
$wsSchedule.Range("H9").Value = 0

# Row 10
$wsSchedule.Range("C10").Value = 42094
$wsSchedule.Range("G10").Value = 3345.38

# Row 11
$wsSchedule.Range("C11").Value = 42095
$wsSchedule.Range("F11").Value = 831.82
$wsSchedule.Range("G11").Value = 2513.56
$wsSchedule.Range("H11").Value = 3.3

# Row 12
$wsSchedule.Range("C12").Value = 42096
$wsSchedule.Range("G12").Value = 1679.27

# Row 13
$wsSchedule.Range("C13").Value = 42097
$wsSchedule.Range("G13").Value = 844.7

# Row 14
$wsSchedule.Range("B14").Value = 1
$wsSchedule.Range("C14").Value = 42098
$wsSchedule.Range("F14").Value = 844.7
$wsSchedule.Range("H14").Value = 0.28000000000000003
$wsSchedule.Range("K14").Value = 844.98
$wsSchedule.Range("Q14").Value = 844.98

# ---------------------------------------------------------------------------
# Sheet: Transactions
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("D1").Value = "TransactionType"
$wsTransactions.Range("A2").Value = 4

# ---------------------------------------------------------------------------
# Selections — restore the recorded cursor position on each sheet. The
# Transactions sheet must stay the active tab, so it is selected last.
# ---------------------------------------------------------------------------
$wsSummary.Range("D7").Select()
$wsSchedule.Range("L9").Select()
$wsTransactions.Range("D7").Select()
